# Applies the cryptos.xlsx update described by the commit:
# "Updated cryptos list on Fri Aug 16 10:47:44 UTC 2024 with GitHub Actions"
#
# Rows 2..51 hold one crypto-coin record each (columns: A=rank, B=Coin,
# C=Link, D=Price, E=Volume(1h)). This run refreshes Price/Volume for
# most rows in place and, for six rows, swaps the B/C/D/E content of two
# adjacent rows (the underlying data source re-ranked those coins) while
# leaving the A (rank index) column untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices as text (e.g. thousand-separated "58.199.25" or
# plain decimals like "6.48"). Excel auto-detects a handful of the new
# prices as real numbers when assigned directly, so those specific cells
# are pre-formatted as Text ("@") to keep them as strings, matching the
# original inline-string cell type.

# --- row 2 ---
$ws.Range("D2").Value = "58.199.25"
$ws.Range("E2").Value = "  -0.13%  "
# --- row 3 ---
$ws.Range("D3").Value = "2.596.24"
$ws.Range("E3").Value = "  -0.53%  "
# --- row 4 ---
$ws.Range("E4").Value = "  +0.15%  "
# --- row 5 ---
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "522.53"
$ws.Range("E5").Value = "  +0.66%  "
# --- row 6 ---
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.57"
$ws.Range("E6").Value = "  +1.43%  "
# --- row 7 ---
$ws.Range("E7").Value = "  -0.19%  "
# --- row 8 ---
$ws.Range("E8").Value = "  +0.78%  "
# --- row 9 ---
$ws.Range("D9").Value = "2.617.06"
$ws.Range("E9").Value = "  -0.08%  "
# --- row 10 ---
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.48"
$ws.Range("E10").Value = "  -0.61%  "
# --- row 11 ---
$ws.Range("E11").Value = "  -0.78%  "
# --- row 12 ---
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.344"
$ws.Range("E12").Value = "  +2.84%  "
# --- row 13 ---
$ws.Range("E13").Value = "  +0.26%  "
# --- row 14 ---
$ws.Range("D14").Value = "3.056.33"
$ws.Range("E14").Value = "  -0.48%  "
# --- row 15 ---
$ws.Range("D15").Value = "58.186.71"
$ws.Range("E15").Value = "  -0.11%  "
# --- row 16 ---
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.40"
$ws.Range("E16").Value = "  -1.97%  "
# --- row 17 ---
$ws.Range("E17").Value = "  -0.73%  "
# --- row 18 ---
$ws.Range("D18").Value = "2.563.50"
$ws.Range("E18").Value = "  -3.41%  "
# --- row 19 ---
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "340.57"
$ws.Range("E19").Value = "  +1.40%  "
# --- row 20 ---
$ws.Range("E20").Value = "  -0.80%  "
# --- row 21 ---
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.25"
$ws.Range("E21").Value = "  -0.89%  "
# --- row 22 ---
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.47"
$ws.Range("E22").Value = "  +3.54%  "
# --- row 23 ---
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  +0.00%  "
# --- row 24 ---
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.44"
$ws.Range("E24").Value = "  +0.98%  "
# --- row 25 ---
$ws.Range("E25").Value = "  +1.89%  "
# --- row 26 ---
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.406"
$ws.Range("E26").Value = "  -1.61%  "
# --- row 27 ---
$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D27").Value = "2.715.69"
$ws.Range("E27").Value = "  -0.66%  "
# --- row 28 ---
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.995"
$ws.Range("E28").Value = "  -0.50%  "
# --- row 29 ---
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.05"
$ws.Range("E29").Value = "  -0.75%  "
# --- row 30 ---
$ws.Range("D30").Value = "0.0₃0746"
$ws.Range("E30").Value = "  -4.58%  "
# --- row 31 ---
$ws.Range("E31").Value = "  -0.11%  "
# --- row 32 ---
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.14"
$ws.Range("E32").Value = "  -5.90%  "
# --- row 33 ---
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.59"
$ws.Range("E33").Value = "  -0.15%  "
# --- row 34 ---
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.82"
$ws.Range("E34").Value = "  +0.66%  "
# --- row 35 ---
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "149.71"
$ws.Range("E35").Value = "  -0.14%  "
# --- row 36 ---
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.01"
$ws.Range("E36").Value = "  -1.20%  "
# --- row 37 ---
$ws.Range("E37").Value = "  -3.42%  "
# --- row 38 ---
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.868"
$ws.Range("E38").Value = "  -2.29%  "
# --- row 39 ---
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.867"
$ws.Range("E39").Value = "  +2.78%  "
# --- row 40 ---
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.47"
$ws.Range("E40").Value = "  +3.19%  "
# --- row 41 ---
$ws.Range("E41").Value = "  -0.32%  "
# --- row 42 ---
$ws.Range("E42").Value = "  -1.66%  "
# --- row 44 ---
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.609"
$ws.Range("E44").Value = "  +1.21%  "
# --- row 45 ---
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "269.87"
$ws.Range("E45").Value = "  +1.22%  "
# --- row 46 ---
$ws.Range("B46").Value = "WhiteBITCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.68"
$ws.Range("E46").Value = "  +0.32%  "
# --- row 47 ---
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0956"
$ws.Range("E47").Value = "  -0.95%  "
# --- row 48 ---
$ws.Range("E48").Value = "  -1.49%  "
# --- row 49 ---
$ws.Range("E49").Value = "  -0.79%  "
# --- row 50 ---
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.77"
$ws.Range("E50").Value = "  +3.65%  "
# --- row 51 ---
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "1.970.03"
$ws.Range("E51").Value = "  -2.50%  "
